# Chapter 1 Notes From Slides
#
# The "Preamble & Scope" bullet under "ABA Model Rules" used to link to the
# in-document bookmark "preamble-scope" (the Heading5 section further down
# the page). It now links out to the ABA's own web page for the Model
# Rules Preamble & Scope, matching the style of its sibling bullets
# (Model Rule 1.2, Model Rule 3.1, Model Rule 3.3), which already point at
# externally hosted americanbar.org pages instead of same-document
# bookmarks.

$d = $word.ActiveDocument

$targetUrl = "http://www.americanbar.org/groups/professional_responsibility/publications/model_rules_of_professional_conduct/model_rules_of_professional_conduct_preamble_scope.html"

foreach ($h in $d.Hyperlinks) {
    if ($h.SubAddress -eq "preamble-scope") {
        $h.Address = $targetUrl
        $h.SubAddress = ""
    }
}
